{"js": "// set_gender(): expand the \"du\u017ean/na je\" placeholder into a Jinja-style\n// {{ du\u017ean_na }} je template tag, matching the gender-suffix convention\n// used elsewhere in the document (e.g. {{ spol_zaposlen_a }}).\n\nconst body = context.document.body;\nconst results = body.search(\"du\u017ean/na je\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n\n  // Build the five-run replacement as raw OOXML so each run keeps its own\n  // run-level formatting, mirroring the diff's run split: a leading\n  // \"apple-converted-space\" run holding \"{{ \", three \"s2\" runs spelling\n  // \"du\u017ean\" + \"_na }}\", and a trailing \"s2\" run holding \" je\". All runs\n  // keep the original Arial/000000/17pt character formatting.\n  const rPrCommon =\n    '<w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n    '<w:color w:val=\"000000\"/><w:sz w:val=\"17\"/><w:szCs w:val=\"17\"/>';\n\n  const run = (style, text) =>\n    `<w:r><w:rPr><w:rStyle w:val=\"${style}\"/>${rPrCommon}</w:rPr>` +\n    `<w:t xml:space=\"preserve\">${text}</w:t></w:r>`;\n\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    run(\"apple-converted-space\", \"{{ \") +\n    run(\"s2\", \"d\") +\n    run(\"s2\", \"u\u017ean\") +\n    run(\"s2\", \"_na }}\") +\n    run(\"s2\", \" je\") +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n\n  // Insert the replacement immediately before the matched range, then\n  // delete the original matched text. Inserting via Word.InsertLocation\n  // .Replace on this range re-anchors at the range's end (leaving stray\n  // merged whitespace behind), so inserting \"before\" and deleting the\n  // original match afterwards keeps everything anchored in place.\n  target.insertOoxml(ooxml, Word.InsertLocation.before);\n  await context.sync();\n\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# set_gender(): expand the \"du\u017ean/na je\" placeholder into a Jinja-style\n# {{ du\u017ean_na }} je template tag, matching the gender-suffix convention\n# used elsewhere in the document (e.g. {{ spol_zaposlen_a }}).\n\n$d = $word.ActiveDocument\n\n$target = $d.Content\n$find = $target.Find\n$find.ClearFormatting()\n$find.Text = \"du\u017ean/na je\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute() | Out-Null\n\nif ($find.Found) {\n    # $target is now collapsed to the matched \"du\u017ean/na je\" text. Clear it,\n    # then rebuild it as five runs mirroring the diff's run-level\n    # formatting split: a leading \"apple-converted-space\" run holding\n    # \"{{ \", three \"s2\" runs spelling \"du\u017ean\" + \"_na }}\", and a trailing\n    # \"s2\" run holding \" je\".\n    $target.Text = \"\"\n\n    $cur = $d.Range($target.Start, $target.Start)\n    $cur.InsertAfter(\"{{ \")\n    $cur.Style = \"apple-converted-space\"\n\n    $cur = $d.Range($cur.End, $cur.End)\n    $cur.InsertAfter(\"d\")\n    $cur.Style = \"s2\"\n\n    $cur = $d.Range($cur.End, $cur.End)\n    $cur.InsertAfter(\"u\u017ean\")\n    $cur.Style = \"s2\"\n\n    $cur = $d.Range($cur.End, $cur.End)\n    $cur.InsertAfter(\"_na }}\")\n    $cur.Style = \"s2\"\n\n    $cur = $d.Range($cur.End, $cur.End)\n    $cur.InsertAfter(\" je\")\n    $cur.Style = \"s2\"\n}\n"}
